$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$findRange = $d.Content
$findRange.Find.Execute("Docente(s) Responsável(eis) ", $true, $false, $false,
                         $false, $false, $true, 1, $false, "", 0)
$headingPara = $findRange.Paragraphs(1)
$headingIndex = $headingPara.Index

# Insert a new paragraph right after it and give it the bullet-list style
# used for similar "name list" entries elsewhere in the document.
$headingPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($headingIndex + 1)
$newPara.Style = "ListBullet"
$newPara.Range.Text = "11079086 - Herlandí de Souza Andrade"
